$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H5").Value = 364.29166
$ws.Range("I5").Value = 170
$ws.Range("K5").Value = 170
$ws.Range("M5").Value = -55
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830

$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 3976.5881
$ws.Range("I45").Value = 5149.5
$ws.Range("K45").Value = 5149.5
$ws.Range("M45").Value = -4772.5
$ws.Range("H76").Value = 29429.334
$ws.Range("J76").Value = 29429.334
$ws.Range("L76").Value = 29429.334
$ws.Range("N76").Value = -30105.334
$ws.Range("H79").Value = 29429.334
$ws.Range("J79").Value = 29429.334
$ws.Range("L79").Value = 29429.334
$ws.Range("N79").Value = -31769.334
$ws.Range("H122").Value = 5673.1
$ws.Range("I122").Value = 3849.625
$ws.Range("K122").Value = 11548.875
$ws.Range("M122").Value = -9098.875
$ws.Range("H132").Value = 7757.891
$ws.Range("I132").Value = 3215.2927
$ws.Range("J132").Value = 45007.2
$ws.Range("K132").Value = 9645.8781
$ws.Range("L132").Value = 135021.6
$ws.Range("M132").Value = -7115.8781
$ws.Range("N132").Value = -140081.6

$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 5881.9614
$ws.Range("I86").Value = 5303.75
$ws.Range("J86").Value = 6807.1
$ws.Range("K86").Value = 5303.75
$ws.Range("L86").Value = 6807.1
$ws.Range("M86").Value = -4180.75
$ws.Range("N86").Value = -9053.1
$ws.Range("H89").Value = 5881.9614
$ws.Range("I89").Value = 5303.75
$ws.Range("J89").Value = 6807.1
$ws.Range("K89").Value = 26518.75
$ws.Range("L89").Value = 34035.5
$ws.Range("M89").Value = -20902.75
$ws.Range("N89").Value = -45267.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 389.64
$ws.Range("I7").Value = 110.82353
$ws.Range("K7").Value = 110.82353
$ws.Range("M7").Value = 2.176469999999995
$ws.Range("H19").Value = 246.42857
$ws.Range("I19").Value = 246.42857
$ws.Range("K19").Value = 246.42857
$ws.Range("M19").Value = -76.42857000000001
$ws.Range("H24").Value = 246.42857
$ws.Range("I24").Value = 246.42857
$ws.Range("K24").Value = 246.42857
$ws.Range("M24").Value = -76.42857000000001
$ws.Range("H31").Value = 28784.682
$ws.Range("I31").Value = 22502.2
$ws.Range("K31").Value = 22502.2
$ws.Range("M31").Value = -22207.2
$ws.Range("H34").Value = 28784.682
$ws.Range("I34").Value = 22502.2
$ws.Range("K34").Value = 22502.2
$ws.Range("M34").Value = -22300.2
$ws.Range("H99").Value = 12421.056
$ws.Range("I99").Value = 6998
$ws.Range("J99").Value = 13505.667
$ws.Range("K99").Value = 6998
$ws.Range("L99").Value = 13505.667
$ws.Range("M99").Value = -5500
$ws.Range("N99").Value = -16501.667
$ws.Range("H109").Value = 18000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 18000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 18000
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -20080
$ws.Range("H126").Value = 12421.056
$ws.Range("I126").Value = 6998
$ws.Range("J126").Value = 13505.667
$ws.Range("K126").Value = 20994
$ws.Range("L126").Value = 40517.001
$ws.Range("M126").Value = -18524
$ws.Range("N126").Value = -45457.001

$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 121.20408
$ws.Range("J2").Value = 61.67857
$ws.Range("L2").Value = 370.07142
$ws.Range("N2").Value = -596.07142
$ws.Range("H11").Value = 2503.6
$ws.Range("I11").Value = 3067.5
$ws.Range("K11").Value = 9202.5
$ws.Range("M11").Value = -9062.5
$ws.Range("H23").Value = 110.6
$ws.Range("J23").Value = 117.666664
$ws.Range("L23").Value = 352.999992
$ws.Range("N23").Value = -822.999992
$ws.Range("H34").Value = 3705571.2
$ws.Range("J34").Value = 11114444
$ws.Range("L34").Value = 33343332
$ws.Range("N34").Value = -33343500
$ws.Range("H39").Value = 7751.2856
$ws.Range("J39").Value = 8432
$ws.Range("L39").Value = 25296
$ws.Range("N39").Value = -25884
$ws.Range("H55").Value = 1734
$ws.Range("J55").Value = 712
$ws.Range("L55").Value = 2136
$ws.Range("N55").Value = -2490
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H126").Value = 7248.6665
$ws.Range("I126").Value = 3164.3333
$ws.Range("J126").Value = 11333
$ws.Range("K126").Value = 9492.999899999999
$ws.Range("L126").Value = 33999
$ws.Range("M126").Value = -4552.999899999999
$ws.Range("N126").Value = -43879
$ws.Range("H130").Value = 13428
$ws.Range("J130").Value = 15332.667
$ws.Range("L130").Value = 45998.001
$ws.Range("N130").Value = -56038.001
$ws.Range("H131").Value = 1421.58
$ws.Range("J131").Value = 1471.914
$ws.Range("L131").Value = 4415.742
$ws.Range("N131").Value = -14495.742
$ws.Range("H138").Value = 4418.84
$ws.Range("I138").Value = 1610
$ws.Range("J138").Value = 4801.864
$ws.Range("K138").Value = 4830
$ws.Range("L138").Value = 14405.592
$ws.Range("M138").Value = 310
$ws.Range("N138").Value = -24685.592

$ws = $wb.Worksheets.Item(6)
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H80").Value = 14190.608
$ws.Range("I80").Value = 11009.5
$ws.Range("J80").Value = 16637.615
$ws.Range("K80").Value = 11009.5
$ws.Range("L80").Value = 16637.615
$ws.Range("M80").Value = -10011.5
$ws.Range("N80").Value = -18633.615
$ws.Range("H83").Value = 14190.608
$ws.Range("I83").Value = 11009.5
$ws.Range("J83").Value = 16637.615
$ws.Range("K83").Value = 55047.5
$ws.Range("L83").Value = 83188.07500000001
$ws.Range("M83").Value = -50055.5
$ws.Range("N83").Value = -93172.07500000001

$ws = $wb.Worksheets.Item(7)
$ws.Range("H20").Value = 261709.81
$ws.Range("I20").Value = 108204
$ws.Range("K20").Value = 108204
$ws.Range("M20").Value = -107978
$ws.Range("H22").Value = 13992.526
$ws.Range("I22").Value = 14639.875
$ws.Range("J22").Value = 13521.728
$ws.Range("K22").Value = 14639.875
$ws.Range("L22").Value = 13521.728
$ws.Range("M22").Value = -14344.875
$ws.Range("N22").Value = -14111.728
$ws.Range("H27").Value = 13992.526
$ws.Range("I27").Value = 14639.875
$ws.Range("J27").Value = 13521.728
$ws.Range("K27").Value = 14639.875
$ws.Range("L27").Value = 13521.728
$ws.Range("M27").Value = -14532.875
$ws.Range("N27").Value = -13735.728
$ws.Range("H46").Value = 2004399.8
$ws.Range("I46").Value = 2503250
$ws.Range("J46").Value = 8999
$ws.Range("K46").Value = 2503250
$ws.Range("L46").Value = 8999
$ws.Range("M46").Value = -2503062
$ws.Range("N46").Value = -9375
$ws.Range("H61").Value = 2470.4285
$ws.Range("I61").Value = 2223.8333
$ws.Range("K61").Value = 2223.8333
$ws.Range("M61").Value = -2021.8333
$ws.Range("H68").Value = 5640.353
$ws.Range("J68").Value = 7719.6
$ws.Range("L68").Value = 7719.6
$ws.Range("N68").Value = -9217.6
$ws.Range("H71").Value = 5640.353
$ws.Range("J71").Value = 7719.6
$ws.Range("L71").Value = 38598
$ws.Range("N71").Value = -46086
$ws.Range("H82").Value = 6765.4443
$ws.Range("I82").Value = 4486
$ws.Range("K82").Value = 4486
$ws.Range("M82").Value = -4125
$ws.Range("H85").Value = 6765.4443
$ws.Range("I85").Value = 4486
$ws.Range("K85").Value = 4486
$ws.Range("M85").Value = -3238
$ws.Range("H113").Value = 2470.4285
$ws.Range("I113").Value = 2223.8333
$ws.Range("K113").Value = 2223.8333
$ws.Range("M113").Value = -53.83329999999978

$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 5250.75
$ws.Range("J62").Value = 5334.3335
$ws.Range("L62").Value = 5334.3335
$ws.Range("N62").Value = -6582.3335
$ws.Range("H65").Value = 5250.75
$ws.Range("J65").Value = 5334.3335
$ws.Range("L65").Value = 26671.6675
$ws.Range("N65").Value = -32911.6675
$ws.Range("H107").Value = 6997.625
$ws.Range("I107").Value = 226.33333
$ws.Range("K107").Value = 678.99999
$ws.Range("M107").Value = 1241.00001
$ws.Range("H122").Value = 4305.1934
$ws.Range("I122").Value = 1761.25
$ws.Range("K122").Value = 5283.75
$ws.Range("M122").Value = -2833.75
$ws.Range("H126").Value = 15428.972
$ws.Range("J126").Value = 11996.167
$ws.Range("L126").Value = 35988.501
$ws.Range("N126").Value = -40928.501
$ws.Range("H132").Value = 8461
$ws.Range("I132").Value = 3657.5
$ws.Range("J132").Value = 22871.5
$ws.Range("K132").Value = 10972.5
$ws.Range("L132").Value = 68614.5
$ws.Range("M132").Value = -8442.5
$ws.Range("N132").Value = -73674.5
$ws.Range("H140").Value = 154702.7
$ws.Range("J140").Value = 154702.7
$ws.Range("L140").Value = 154702.7
$ws.Range("N140").Value = -165062.7
